$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.693909287452698
$ws.Range("B1").Value = 3.757834196090698
$ws.Range("C1").Value = 2.967264652252197
$ws.Range("D1").Value = 2.396448373794556
$ws.Range("E1").Value = 1.448315620422363
